$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 6
$ws.Range("H6").Value = 100
$ws.Range("I6").Value = 100
$ws.Range("K6").Value = 300
$ws.Range("M6").Value = -188

# ALC row 70
$ws.Range("H70").Value = 7377.1665
$ws.Range("J70").Value = 8248.5
$ws.Range("L70").Value = 24745.5
$ws.Range("N70").Value = -25285.5

# ALC row 73
$ws.Range("H73").Value = 7377.1665
$ws.Range("J73").Value = 8248.5
$ws.Range("L73").Value = 24745.5
$ws.Range("N73").Value = -26617.5

# ALC row 94
$ws.Range("H94").Value = 14997.5
$ws.Range("I94").Value = 14997.5
$ws.Range("K94").Value = 14997.5
$ws.Range("M94").Value = -14546.5

# ALC row 125
$ws.Range("H125").Value = 3055.2307
$ws.Range("J125").Value = 3923.1428
$ws.Range("L125").Value = 35308.2852
$ws.Range("N125").Value = -40228.2852

# ALC row 132
$ws.Range("H132").Value = 1358.4595
$ws.Range("I132").Value = 1213.9412
$ws.Range("K132").Value = 3641.8236
$ws.Range("M132").Value = -1111.8236

# ALC row 137
$ws.Range("H137").Value = 6952287
$ws.Range("I137").Value = 20838250
$ws.Range("J137").Value = 9304.916999999999
$ws.Range("K137").Value = 62514750
$ws.Range("L137").Value = 27914.751
$ws.Range("M137").Value = -62512200
$ws.Range("N137").Value = -33014.751

# ALC row 138
$ws.Range("H138").Value = 4556.0537
$ws.Range("I138").Value = 4211.0625
$ws.Range("J138").Value = 4694.05
$ws.Range("K138").Value = 12633.1875
$ws.Range("L138").Value = 14082.15
$ws.Range("M138").Value = -7493.1875
$ws.Range("N138").Value = -24362.15

$ws = $wb.Worksheets.Item("BSM")
# BSM row 5
$ws.Range("H5").Value = 5000
$ws.Range("I5").Value = 5000
$ws.Range("J5").Value = 5000
$ws.Range("K5").Value = 5000
$ws.Range("L5").Value = 5000
$ws.Range("M5").Value = -4887
$ws.Range("N5").Value = -5226

# BSM row 54
$ws.Range("H54").Value = 13079.833
$ws.Range("I54").Value = 13079.833
$ws.Range("K54").Value = 13079.833
$ws.Range("M54").Value = -12595.833

# BSM row 134
$ws.Range("H134").Value = 4609.2
$ws.Range("I134").Value = 2352.0645
$ws.Range("J134").Value = 9607.143
$ws.Range("K134").Value = 7056.193499999999
$ws.Range("L134").Value = 28821.429
$ws.Range("M134").Value = -4521.193499999999
$ws.Range("N134").Value = -33891.429

$ws = $wb.Worksheets.Item("CRP")
# CRP row 12
$ws.Range("H12").Value = 120999.8
$ws.Range("I12").Value = 120999.8
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 120999.8
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -120829.8
$ws.Range("N12").ClearContents()

# CRP row 16
$ws.Range("H16").Value = 2305.0588
$ws.Range("I16").Value = 1699.6
$ws.Range("J16").Value = 2783.0527
$ws.Range("K16").Value = 1699.6
$ws.Range("L16").Value = 2783.0527
$ws.Range("M16").Value = -1412.6
$ws.Range("N16").Value = -3357.0527

# CRP row 31
$ws.Range("H31").Value = 31253216
$ws.Range("I31").Value = 83335460
$ws.Range("K31").Value = 83335460
$ws.Range("M31").Value = -83335165

# CRP row 34
$ws.Range("H34").Value = 31253216
$ws.Range("I34").Value = 83335460
$ws.Range("K34").Value = 83335460
$ws.Range("M34").Value = -83335258

# CRP row 62
$ws.Range("H62").Value = 19998.54
$ws.Range("J62").Value = 20888.223
$ws.Range("L62").Value = 20888.223
$ws.Range("N62").Value = -22136.223

# CRP row 65
$ws.Range("H65").Value = 19998.54
$ws.Range("J65").Value = 20888.223
$ws.Range("L65").Value = 104441.115
$ws.Range("N65").Value = -110681.115

# CRP row 86
$ws.Range("H86").Value = 6505.636
$ws.Range("I86").Value = 5392.6665
$ws.Range("J86").Value = 7841.2
$ws.Range("K86").Value = 5392.6665
$ws.Range("L86").Value = 7841.2
$ws.Range("M86").Value = -4269.6665
$ws.Range("N86").Value = -10087.2

# CRP row 89
$ws.Range("H89").Value = 6505.636
$ws.Range("I89").Value = 5392.6665
$ws.Range("J89").Value = 7841.2
$ws.Range("K89").Value = 26963.3325
$ws.Range("L89").Value = 39206
$ws.Range("M89").Value = -21347.3325
$ws.Range("N89").Value = -50438

# CRP row 93
$ws.Range("H93").Value = 31914.666
$ws.Range("I93").Value = 16872
$ws.Range("J93").Value = 62000
$ws.Range("K93").Value = 16872
$ws.Range("L93").Value = 62000
$ws.Range("M93").Value = -15000
$ws.Range("N93").Value = -65744

# CRP row 113
$ws.Range("H113").Value = 2305.0588
$ws.Range("I113").Value = 1699.6
$ws.Range("J113").Value = 2783.0527
$ws.Range("K113").Value = 1699.6
$ws.Range("L113").Value = 2783.0527
$ws.Range("M113").Value = 470.4000000000001
$ws.Range("N113").Value = -7123.0527

$ws = $wb.Worksheets.Item("CUL")
# CUL row 44
$ws.Range("H44").Value = 1912.909
$ws.Range("I44").Value = 799.75
$ws.Range("J44").Value = 2549
$ws.Range("K44").Value = 2399.25
$ws.Range("L44").Value = 7647
$ws.Range("M44").Value = -2001.25
$ws.Range("N44").Value = -8443

# CUL row 68
$ws.Range("H68").Value = 56117.05
$ws.Range("I68").Value = 287885.44
$ws.Range("J68").Value = 3782.258
$ws.Range("K68").Value = 863656.3200000001
$ws.Range("L68").Value = 11346.774
$ws.Range("M68").Value = -862845.3200000001
$ws.Range("N68").Value = -12968.774

# CUL row 71
$ws.Range("H71").Value = 56117.05
$ws.Range("I71").Value = 287885.44
$ws.Range("J71").Value = 3782.258
$ws.Range("K71").Value = 2590968.96
$ws.Range("L71").Value = 34040.322
$ws.Range("M71").Value = -2586912.96
$ws.Range("N71").Value = -42152.322

# CUL row 113
$ws.Range("H113").Value = 2693.8823
$ws.Range("I113").Value = 2060.5
$ws.Range("J113").Value = 2888.7693
$ws.Range("K113").Value = 6181.5
$ws.Range("L113").Value = 8666.3079
$ws.Range("M113").Value = -4011.5
$ws.Range("N113").Value = -13006.3079

# CUL row 132
$ws.Range("H132").Value = 1879.1111
$ws.Range("I132").Value = 1513.7142
$ws.Range("J132").Value = 2111.6365
$ws.Range("K132").Value = 13623.4278
$ws.Range("L132").Value = 19004.7285
$ws.Range("M132").Value = -11093.4278
$ws.Range("N132").Value = -24064.7285

# CUL row 137
$ws.Range("H137").Value = 4056.4614
$ws.Range("J137").Value = 4903.8887
$ws.Range("L137").Value = 14711.6661
$ws.Range("N137").Value = -24911.6661

$ws = $wb.Worksheets.Item("GSM")
# GSM row 113
$ws.Range("H113").Value = 19832.834
$ws.Range("I113").Value = 3427.7144
$ws.Range("J113").Value = 42800
$ws.Range("K113").Value = 3427.7144
$ws.Range("L113").Value = 42800
$ws.Range("M113").Value = -1257.7144
$ws.Range("N113").Value = -47140

# GSM row 126
$ws.Range("H126").Value = 3077.9583
$ws.Range("I126").Value = 2553.0588
$ws.Range("K126").Value = 7659.176399999999
$ws.Range("M126").Value = -5189.176399999999

$ws = $wb.Worksheets.Item("LTW")
# LTW row 68
$ws.Range("H68").Value = 2223.75
$ws.Range("I68").Value = 2223.75
$ws.Range("K68").Value = 2223.75
$ws.Range("M68").Value = -1474.75

# LTW row 71
$ws.Range("H71").Value = 2223.75
$ws.Range("I71").Value = 2223.75
$ws.Range("K71").Value = 11118.75
$ws.Range("M71").Value = -7374.75

$ws = $wb.Worksheets.Item("WVR")
# WVR row 17
$ws.Range("H17").Value = 100500
$ws.Range("I17").Value = 1000
$ws.Range("K17").Value = 1000
$ws.Range("M17").Value = -828

# WVR row 51
$ws.Range("H51").Value = 168852.33
$ws.Range("I51").Value = 209045.19
$ws.Range("K51").Value = 209045.19
$ws.Range("M51").Value = -208535.19

# WVR row 75
$ws.Range("H75").Value = 70000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 70000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 70000
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -71872

# WVR row 78
$ws.Range("H78").Value = 70000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 70000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 210000
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -219360

# WVR row 95
$ws.Range("H95").Value = 45855.145
$ws.Range("J95").Value = 45855.145
$ws.Range("L95").Value = 45855.145
$ws.Range("N95").Value = -51347.145

# WVR row 136
$ws.Range("H136").Value = 5618.8
$ws.Range("I136").Value = 3914.7778
$ws.Range("J136").Value = 8174.8335
$ws.Range("K136").Value = 11744.3334
$ws.Range("L136").Value = 24524.5005
$ws.Range("M136").Value = -9194.3334
$ws.Range("N136").Value = -29624.5005
